# add lista de pedidos PDF
# Mark REQ-037, REQ-038, REQ-039 (rows 9, 10, 11) as collected/concluded via
# the new automated "PDF Exportado" flow: update status, last-update
# timestamp, and fill in the separation/collection responsible + dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (REQ-037): was PENDENTE, now CONCLUÍDO via PDF export ---
$ws.Range("K9").Value = "CONCLUÍDO"
$ws.Range("M9").Value = "02/07/2025 14:28"
$ws.Range("N9").Value = "Sistema (PDF Exportado)"
$ws.Range("O9").Value = "Sistema (PDF Exportado)"
$ws.Range("P9").Value = "02/07/2025 14:12"
$ws.Range("Q9").Value = "Sistema (PDF Exportado)"
$ws.Range("R9").Value = "02/07/2025 14:28"

# --- Row 10 (REQ-038): already CONCLUÍDO, collection now via PDF export ---
$ws.Range("M10").Value = "02/07/2025 14:28"
$ws.Range("N10").Value = "Sistema (PDF Exportado)"
$ws.Range("O10").Value = "Sistema (PDF Exportado)"
$ws.Range("P10").Value = "02/07/2025 14:13"
$ws.Range("Q10").Value = "Sistema (PDF Exportado)"
$ws.Range("R10").Value = "02/07/2025 14:28"

# --- Row 11 (REQ-039): already CONCLUÍDO, collection now via PDF export ---
$ws.Range("M11").Value = "02/07/2025 14:28"
$ws.Range("N11").Value = "Sistema (PDF Exportado)"
$ws.Range("O11").Value = "Sistema (PDF Exportado)"
$ws.Range("P11").Value = "02/07/2025 14:13"
$ws.Range("Q11").Value = "Sistema (PDF Exportado)"
$ws.Range("R11").Value = "02/07/2025 14:28"
